# The two observation records that used to be on rows 2 and 3 were
# reordered: what used to be row 2 is now row 3, and vice-versa.
# Row 1 (headers) and row 4 (the third observation) are untouched.
#
# We swap cell-by-cell across the full used width of the sheet, but we
# only actually write to a cell when its value truly changes between the
# two rows. This avoids pointless writes to already-identical cells
# (e.g. blank cells, or text that Excel could otherwise try to
# reinterpret, such as "2023-08-31" date strings that are identical on
# both rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2Num = 2
$row3Num = 3
$lastCol = $ws.UsedRange.Columns.Count
if ($lastCol -lt 51) { $lastCol = 51 }   # sheet data goes through column AY

for ($c = 1; $c -le $lastCol; $c++) {
    $cell2 = $ws.Cells.Item($row2Num, $c)
    $cell3 = $ws.Cells.Item($row3Num, $c)

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    if ($v2 -ne $v3) {
        $cell2.Value = $v3
        $cell3.Value = $v2
    }
}
